# Auto-generated edit script: apply "Update latest output (run 134)" changes
$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

$wsSchedule.Range("A3").Value = 46046.25
$wsSchedule.Range("B3").Value = 46046.75
$wsSchedule.Range("E3").Value = -31.60218749999996
$wsSchedule.Range("F3").Value = -0.6966972552910045
$wsSchedule.Range("E4").Value = 363.4444515000001
$wsSchedule.Range("F4").Value = 8.012443816137568
$wsDetailed.Range("E14").Value = "ON"
$wsDetailed.Range("B37").Value = 49.24117
$wsDetailed.Range("B38").Value = 66.35997
$wsDetailed.Range("E38").Value = "OFF"
$wsDetailed.Range("B39").Value = 83.01918999999999
$wsDetailed.Range("C39").Value = "historical"
$wsDetailed.Range("B40").Value = 119.83444
$wsDetailed.Range("C40").Value = "historical"
$wsDetailed.Range("B41").Value = 147.04151
$wsDetailed.Range("C41").Value = "historical"
$wsDetailed.Range("B42").Value = 147.51746
$wsDetailed.Range("C42").Value = "historical"
$wsDetailed.Range("B43").Value = 186.23306
$wsDetailed.Range("C43").Value = "historical"
$wsDetailed.Range("B44").Value = 138.42
$wsDetailed.Range("C44").Value = "historical"
$wsDetailed.Range("B45").Value = 116.49262
$wsDetailed.Range("C45").Value = "historical"
$wsDetailed.Range("B46").Value = 108.89
$wsDetailed.Range("C46").Value = "historical"
$wsDetailed.Range("B47").Value = 101.25
$wsDetailed.Range("C47").Value = "historical"
$wsDetailed.Range("B48").Value = 102.30723
$wsDetailed.Range("C48").Value = "historical"
$wsDetailed.Range("B49").Value = 105.79
$wsDetailed.Range("B50").Value = 105.79
$wsDetailed.Range("B51").Value = 124.32627
$wsDetailed.Range("B52").Value = 105.79
$wsDetailed.Range("B53").Value = 105.79
$wsDetailed.Range("B54").Value = 103.78068
$wsDetailed.Range("B55").Value = 105.00005
$wsDetailed.Range("B56").Value = 105.00005
$wsDetailed.Range("B57").Value = 105.35984
$wsDetailed.Range("B58").Value = 87.60003
$wsDetailed.Range("B59").Value = 98.25972
$wsDetailed.Range("B60").Value = 103.39668
$wsDetailed.Range("B61").Value = 102.36287
$wsDetailed.Range("B62").Value = 99.50184
$wsDetailed.Range("B63").Value = 62.29726
$wsDetailed.Range("B64").Value = 51.22225
$wsDetailed.Range("B65").Value = 19.80149
$wsDetailed.Range("B66").Value = 5.70924
$wsDetailed.Range("B67").Value = -1.06044
$wsDetailed.Range("B68").Value = -5.50985
$wsDetailed.Range("B69").Value = -6
$wsDetailed.Range("B70").Value = -6.19359
$wsDetailed.Range("B71").Value = -6.16627
$wsDetailed.Range("B72").Value = -6.08061
$wsDetailed.Range("B73").Value = -6
$wsDetailed.Range("B74").Value = -6
$wsDetailed.Range("B75").Value = -6.52693
$wsDetailed.Range("B76").Value = -6
$wsDetailed.Range("B77").Value = -5.50985
$wsDetailed.Range("B78").Value = 0.00002
$wsDetailed.Range("B79").Value = -2.54301
$wsDetailed.Range("B80").Value = -1.49771
$wsDetailed.Range("B81").Value = 0.41892
$wsDetailed.Range("B82").Value = 2.68504
$wsDetailed.Range("B83").Value = 48.75908
$wsDetailed.Range("B84").Value = 70.50918
$wsDetailed.Range("B85").Value = 76.94748
$wsDetailed.Range("B86").Value = 184.43252
$wsDetailed.Range("B87").Value = 349.81754
$wsDetailed.Range("B88").Value = 408.81344
$wsDetailed.Range("B89").Value = 286.65072
$wsDetailed.Range("B90").Value = 330.17155
$wsDetailed.Range("B91").Value = 168.29492
$wsDetailed.Range("B92").Value = 279.10835
$wsDetailed.Range("B93").Value = 240.89
$wsDetailed.Range("B94").Value = 138.57805
$wsDetailed.Range("B95").Value = 147.52
$wsDetailed.Range("B97").Value = 91.65622999999999
